$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 — copy formatting (bold, border, centered/top alignment)
# from the existing header cell E1 so it reuses the same cell style, then
# set its value.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# New data column values (plain, unstyled like the rest of column F/E body cells)
$ws.Range("F2").Value = "2021-10-05 13:41:47.345750"
$ws.Range("F3").Value = "2021-10-05 13:41:47.345763"
$ws.Range("F4").Value = "2021-10-05 13:41:47.345766"
$ws.Range("F5").Value = "2021-10-05 13:41:47.345769"
